$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column header "avg%" (C1) - reuse same centered style as existing headers
$ws.Range("C1").Value = "avg%"
$ws.Range("C1").HorizontalAlignment = -4108

# C2 gets its own (non-shared) formula
$ws.Range("C2").Formula = "=ABS(B2/1000)"

# C3:C8 share one formula (anchored at C3) mirroring B3:B8
$ws.Range("C3:C8").Formula = "=ABS(B3/1000)"

# Apply the same centered style to C2:C8
$ws.Range("C2:C8").HorizontalAlignment = -4108

# B9 becomes a styled-but-empty cell
$ws.Range("B9").HorizontalAlignment = -4108

# Column B narrows from 34.140625 to 31 (stored units); 30.1666... chars maps to 31
$ws.Columns("B").ColumnWidth = 30.166666666666668

# Move the active selection to F3
$ws.Range("F3").Select()
